# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps recorded for the zh-cn and de-de handback rows.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 22:51:18"
$wsZhCn.Range("H2").Value = "2016-03-19 22:51:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 22:51:21"
$wsDeDe.Range("H2").Value = "2016-03-19 22:51:42"
